$d = $word.ActiveDocument

$pairs = @(
    @("81×25=", "31×19="),
    @("60×27=", "24×19="),
    @("94×21=", "31×17="),
    @("71×59=", "71×99="),
    @("87×23=", "88×24="),
    @("11×13=", "20×62="),
    @("63×11=", "22×45="),
    @("92×79=", "72×19="),
    @("94×62=", "70×58="),
    @("87×32=", "60×12="),
    @("66×47=", "88×99="),
    @("94×34=", "25×69="),
    @("36×98=", "37×37="),
    @("36×95=", "29×66="),
    @("23×39=", "77×62="),
    @("63×16=", "21×67="),
    @("77×77=", "33×72="),
    @("67×91=", "82×30="),
    @("60×65=", "59×34="),
    @("46×57=", "95×83="),
    @("90×96=", "38×81="),
    @("63×21=", "18×11="),
    @("71×74=", "59×79="),
    @("21×61=", "92×12="),
    @("72×21=", "30×55=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
